# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Both sheets contain identical data, so the same row -> value updates apply to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3080
    5  = 162
    7  = 1700
    12 = 1393
    14 = 534
    15 = 351
    16 = 49
    17 = 8
    21 = 92
    22 = 108
    23 = 3259
    24 = 396
    25 = 159
    26 = 337
    28 = 23
    29 = 117
    30 = 101
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

$wb.Save()
